# The document's last paragraph currently reads:
#   "多云转小雨，今天学习了多分支管理，创建了一个dev分支。"
# We split it into two paragraphs:
#   1) a new paragraph carrying the original sentence verbatim
#   2) the original (last) paragraph, now holding the new sentence
#      "使用dev创建分支简单又快速。" -- keeping its bookmarkStart/bookmarkEnd.

$d = $word.ActiveDocument

$original = "多云转小雨，今天学习了多分支管理，创建了一个dev分支。"
$addition = "使用dev创建分支简单又快速。"

# Locate the paragraph that currently holds the original sentence (the
# document's final paragraph) and the paragraph right before it, whose
# paragraph-mark formatting (rFonts hint="eastAsia") the newly split-off
# paragraph should inherit.
$lastIndex = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($lastIndex - 1)
$prevRange = $prevPara.Range

# Insert a paragraph break right after the previous paragraph's text (but
# before its own paragraph mark), which mints a new empty paragraph that
# inherits that paragraph's "eastAsia" hinted formatting -- matching how
# the target XML formats the newly split-out paragraph.
$splitPoint = $d.Range($prevRange.End - 1, $prevRange.End - 1)
$splitPoint.InsertParagraphAfter()

# Fill the freshly minted paragraph with the original sentence.
$newPara = $d.Paragraphs.Item($lastIndex)
$newPara.Range.Text = $original

# The original last paragraph (now shifted one index later) still has the
# old sentence plus the _GoBack bookmark; replace its text in place so the
# bookmark markers stay attached to this paragraph.
$finalPara = $d.Paragraphs.Item($lastIndex + 1)
$finalPara.Range.Find.Execute($original, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $addition, 2)
